# Harmonize similar tags to be the same.
#
# The "isa_template" sheet holds the Swate template metadata. The #TAGS
# list row ("Tags" / row 12 together with its accession-number / term-
# source rows 13-14) listed four separate tag fragments:
#   B12 = "Proteomics"
#   C12 = "Sample"
#   D12 = "Mass spectrometry"
#   E12 = "MS"
#   F12 = "Preparation"
# "Mass spectrometry" / "MS" and "Sample" / "Preparation" were really the
# same two tags, just written inconsistently (capitalisation / abbreviated
# form) and split across an extra column. Harmonize them: keep "Sample" in
# C12, fix "Mass spectrometry" to "Mass Spectrometry" in D12, fold the old
# F12 ("Preparation") value into E12 (replacing the stray "MS"), and drop
# the now redundant F12 cell. The matching accession-number row (13) gets
# the proper NCIT accession for the harmonized "Mass Spectrometry" tag.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_template")

# Row 12 ("Tags"): fix casing of the "Mass spectrometry" tag and replace
# the stray "MS" abbreviation with the "Preparation" tag that used to live
# in the now-removed F12 cell.
$ws.Range("D12").Value = "Mass Spectrometry"
$ws.Range("E12").Value = "Preparation"
$ws.Range("F12").Clear()

# Row 13 ("Tags Term Accession Number"): record the accession number for
# the harmonized "Mass Spectrometry" tag.
$ws.Range("D13").Value = "NCIT:C17156"

# Row heights for rows 12/13 change now that the content reflows (one
# fewer wrapped column in row 12, new wrapped text in row 13).
$ws.Rows.Item(12).RowHeight = 43.2
$ws.Rows.Item(13).RowHeight = 28.8

# Reflect the last-used selection on the "isa_template" sheet.
$ws.Activate() | Out-Null
$ws.Range("F9").Select() | Out-Null
